# Weekly price-list update: a new daily record is inserted at row 171
# (pushing the existing rows 171-261 down to 172-262, and growing the
# sheet's used range from A1:R261 to A1:R262).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 171..261 down by one to make room for the new record.
$ws.Rows.Item(171).Insert()

# Populate the newly inserted row 171 with the new observation.
$ws.Cells.Item(171, 1).Value = 7
$ws.Cells.Item(171, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(171, 3).Value = "Ñuble"
$ws.Cells.Item(171, 4).Value = 44839
$ws.Cells.Item(171, 5).Value = 16
$ws.Cells.Item(171, 6).Value = 100112043
$ws.Cells.Item(171, 7).Value = "Pepino ensalada"
$ws.Cells.Item(171, 8).Value = "Sin especificar"
$ws.Cells.Item(171, 9).Value = "Primera"
$ws.Cells.Item(171, 10).Value = 50
$ws.Cells.Item(171, 11).Value = 25000
$ws.Cells.Item(171, 12).Value = 25000
$ws.Cells.Item(171, 13).Value = 25000
$ws.Cells.Item(171, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(171, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(171, 16).Value = 417
$ws.Cells.Item(171, 17).Value = 60
$ws.Cells.Item(171, 18).Value = "Hortaliza"

# Make sure the date column keeps the workbook's date number format
# (it is normally inherited from the row above on insert, but set it
# explicitly to be safe).
$ws.Cells.Item(171, 4).NumberFormat = $ws.Cells.Item(172, 4).NumberFormat
